$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The organization's website changed from www.stat.kg to www.stat.gov.kg.
# (Row 10, "Сайт организации (если есть)")
$ws.Range("B10").Value = "www.stat.gov.kg"

# Re-apply direct font formatting on the three cells that were touched
# during this editing pass (B4, B9, B10) - this is what produced the
# extra, slightly-different cell style ("s=7") Excel wrote for them on
# save, distinct from the untouched cells that still share style "s=5".
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B10").Font.Name = "Calibri"

# Leave the cursor on the last-edited cell, matching the saved selection.
$ws.Range("B4").Select()
